$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above the current row 1 so the existing header
# row (row 1) and the two data rows (rows 2-3) shift down to rows 4-6.
$ws.Range("A1:A3").EntireRow.Insert()

# New row 1: warehouse label
$ws.Range("A1").Value = "Sklad:"

# New row 2: warehouse name
$ws.Range("A2").Value = "prvy"

# Row 3 stays empty (gap row).

# Row 4 (formerly row 1): translated column headers
$ws.Range("A4").Value = "ID produktu"
$ws.Range("B4").Value = "Meno produktu"
$ws.Range("C4").Value = "Pocet"
$ws.Range("D4").Value = "Cena za jednotku"
$ws.Range("E4").Value = "Jednotková váha"
$ws.Range("F4").Value = "Vlastnosti"
$ws.Range("G4").Value = "Link na obrázok"
